# #clean code and fix error
# Remove stale "&=result.X" placeholder formulas for columns that are no
# longer populated (Status, Hot_Sale, New, IsSale) and clear the
# corresponding cells in row 2 (columns C:F), keeping only the
# Product Cate ID (A2) and Product Name (B2) placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused placeholder cells (C2:F2) while leaving their
# style/borders intact.
$ws.Range("C2:F2").ClearContents()

# Update the active selection to reflect where the cursor ended up.
$ws.Range("F2").Select()
